# -----------------------------------------------------------------------
# Applies the "optimization_parameters" sheet addition to the
# expression_sheet_empty_column workbook, as described by the commit:
#   "edited all expression sheet test files to include species name and
#    taxon id"
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add a new worksheet "optimization_parameters" as the last sheet.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "optimization_parameters"

# ---------------------------------------------------------------------
# 2. Populate the data for the new sheet, row by row. Each inner array
#    holds the values for one row, starting in column A.
# ---------------------------------------------------------------------
$rows = @(
    @("optimization_parameter", "value"),
    @("alpha", 0.002),
    @("kk_max", 1),
    @("MaxIter", 100000000),
    @("TolFun", 0.000001),
    @("MaxFunEval", 100000000),
    @("TolX", 0.000001),
    @("production_function", "Sigmoid"),
    @("L_curve", 0),
    @("estimate_params", 1),
    @("make_graphs", 1),
    @("fix_P", 0),
    @("fix_b", 0),
    @("expression_timepoints", 15, 30, 60),
    @("Strain", "wt", "dgln3"),
    @("simulation_timepoints", 0, 5, 10, 15),
    @("species", "Saccharomyces cerevisiae"),
    @("taxon_id", 559292)
)

$rowIndex = 1
foreach ($rowValues in $rows) {
    $colIndex = 1
    foreach ($val in $rowValues) {
        $ws.Cells.Item($rowIndex, $colIndex).Value = $val
        $colIndex = $colIndex + 1
    }

    # Format just the populated cells of this row (Arial 10pt black),
    # leaving the rest of the row untouched / blank.
    $rowRange = $ws.Range($ws.Cells.Item($rowIndex, 1), $ws.Cells.Item($rowIndex, $colIndex - 1))
    $rowRange.Font.Name = "Arial"
    $rowRange.Font.Size = 10
    $rowRange.Font.Color = 0

    $rowIndex = $rowIndex + 1
}

# ---------------------------------------------------------------------
# 3. The five numeric solver-tolerance parameters are displayed using
#    scientific notation.
# ---------------------------------------------------------------------
foreach ($addr in @("B2", "B4", "B5", "B6", "B7")) {
    $ws.Range($addr).NumberFormat = "0.00E+00"
}

# Column A is wider to fit the parameter names.
$ws.Columns.Item(1).ColumnWidth = 21.5

# ---------------------------------------------------------------------
# 4. Make the new sheet the active / selected sheet, mirroring the
#    workbook's stored selection in the diff.
# ---------------------------------------------------------------------
$ws.Activate()
[void]$ws.Range("C22").Select()

# ---------------------------------------------------------------------
# 5. Workbook-level iterative-calculation setting referenced by the
#    diff (calcPr iterateDelta="1E-4").
# ---------------------------------------------------------------------
$excel.Iteration = $true
$excel.MaxChange = 0.0001

$wb.Save()
